$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# --- Cells that flip from a numeric value to the "0"/"***.*" placeholder text ---
# These reuse the exact cell style already used for similar placeholder cells
# (C14 carries the "0" text style, E14 carries the "***.*" text style) by
# writing the text first (quote-prefixed so Excel keeps it as text instead of
# re-parsing "0" back into a number), then copying *only* the formatting from
# the reference cell so the resulting style index matches exactly.

function Set-PlaceholderText {
    param([string]$TargetCell, [string]$Text, [string]$StyleSourceCell)
    $ws.Range($TargetCell).Value = "'" + $Text
    $ws.Range($StyleSourceCell).Copy() | Out-Null
    $ws.Range($TargetCell).PasteSpecial(-4122) | Out-Null
}

Set-PlaceholderText "G14" "0" "C14"
Set-PlaceholderText "H14" "***.*" "E14"
Set-PlaceholderText "D15" "0" "C14"
Set-PlaceholderText "E15" "***.*" "E14"
Set-PlaceholderText "C22" "0" "C14"

# --- Remaining plain numeric updates ---
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -58.474576271186
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 56
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = 9.803921568627
$ws.Range("I16").Value = 712
$ws.Range("J16").Value = 579
$ws.Range("K16").Value = 22.970639032815
$ws.Range("L16").Value = 25.352112676056
$ws.Range("M16").Value = 1.569186875891
$ws.Range("N16").Value = -75.482093663911
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -11.764705882352
$ws.Range("F17").Value = 74
$ws.Range("G17").Value = 89
$ws.Range("H17").Value = -16.853932584269
$ws.Range("I17").Value = 1018
$ws.Range("J17").Value = 888
$ws.Range("K17").Value = 14.639639639639
$ws.Range("L17").Value = 24.907975460122
$ws.Range("M17").Value = 44.397163120567
$ws.Range("N17").Value = -26.709863210943
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = -44.736842105263
$ws.Range("I18").Value = 353
$ws.Range("J18").Value = 380
$ws.Range("K18").Value = -7.105263157894
$ws.Range("L18").Value = -14.320388349514
$ws.Range("M18").Value = -21.902654867256
$ws.Range("N18").Value = -79.712643678160
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -48.484848484848
$ws.Range("F19").Value = 81
$ws.Range("G19").Value = 111
$ws.Range("H19").Value = -27.027027027027
$ws.Range("I19").Value = 1226
$ws.Range("J19").Value = 997
$ws.Range("K19").Value = 22.968906720160
$ws.Range("L19").Value = 36.222222222222
$ws.Range("M19").Value = 85.196374622356
$ws.Range("N19").Value = 37.443946188340
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 46
$ws.Range("H20").Value = 31.428571428571
$ws.Range("I20").Value = 482
$ws.Range("J20").Value = 492
$ws.Range("K20").Value = -2.032520325203
$ws.Range("L20").Value = 30.270270270270
$ws.Range("M20").Value = 67.361111111111
$ws.Range("N20").Value = -80.350591112923
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 85
$ws.Range("E21").Value = -30.588235294117
$ws.Range("F21").Value = 284
$ws.Range("G21").Value = 330
$ws.Range("H21").Value = -13.939393939393
$ws.Range("I21").Value = 3863
$ws.Range("J21").Value = 3414
$ws.Range("K21").Value = 13.151728178090
$ws.Range("L21").Value = 22.557106598984
$ws.Range("M21").Value = 33.344839489126
$ws.Range("N21").Value = -59.785550697480
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = -27.272727272727
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = 52.083333333333
$ws.Range("L22").Value = 19.672131147541
$ws.Range("M22").Value = 12.307692307692
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = -54.545454545454
$ws.Range("F23").Value = 29
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = 3.571428571428
$ws.Range("I23").Value = 331
$ws.Range("J23").Value = 283
$ws.Range("K23").Value = 16.961130742049
$ws.Range("L23").Value = 32.4
$ws.Range("M23").Value = 53.953488372093
$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = 23.076923076923
$ws.Range("F24").Value = 223
$ws.Range("G24").Value = 202
$ws.Range("H24").Value = 10.396039603960
$ws.Range("I24").Value = 2714
$ws.Range("J24").Value = 2167
$ws.Range("K24").Value = 25.242270419935
$ws.Range("L24").Value = 45.211342964152
$ws.Range("M24").Value = 95.533141210374
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 96
$ws.Range("H25").Value = -14.583333333333
$ws.Range("I25").Value = 1183
$ws.Range("J25").Value = 972
$ws.Range("K25").Value = 21.707818930041
$ws.Range("L25").Value = 39.504716981132
$ws.Range("M25").Value = -31.300813008130
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = -50
$ws.Range("I26").Value = 75
$ws.Range("J26").Value = 94
$ws.Range("K26").Value = -20.212765957446
$ws.Range("L26").Value = -14.772727272727
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -87.5
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = -43.75
$ws.Range("I27").Value = 112
$ws.Range("J27").Value = 119
$ws.Range("K27").Value = -5.882352941176
$ws.Range("L27").Value = 13.131313131313
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 88
$ws.Range("J28").Value = 81
$ws.Range("K28").Value = 8.641975308641
$ws.Range("L28").Value = -24.137931034482
$ws.Range("M28").Value = -16.981132075471
$ws.Range("N28").Value = -69.444444444444
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -66.666666666666
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 33.333333333333
$ws.Range("I29").Value = 77
$ws.Range("J29").Value = 67
$ws.Range("K29").Value = 14.925373134328
$ws.Range("L29").Value = -18.947368421052
$ws.Range("M29").Value = -11.494252873563
$ws.Range("N29").Value = -70.498084291187
